# Before each test method, the driver is wired up with a different proxy
# IP address. Rotate the "LoginData" sheet's proxy list (column B, rows
# 2-11) to a new set of addresses for this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")

$ws.Range("B2").Value = "185.157.1.249:3128"
$ws.Range("B4").Value = "45.169.162.1:3128"
$ws.Range("B5").Value = "157.100.12.138:999"
$ws.Range("B6").Value = "157.230.217.232:8080"
$ws.Range("B7").Value = "177.12.238.1:3128"
$ws.Range("B8").Value = "116.58.254.126:8080"
$ws.Range("B9").Value = "89.107.197.165:3128"
$ws.Range("B10").Value = "191.97.19.18:999"
$ws.Range("B11").Value = "177.12.238.100:3128"
